# "some excel import fixes"
# - B2 on Sheet1 changes from the numeric value 11001028583 to the text
#   value "123456789" (stored as a new shared string).
# - The sheet's active selection moves from L10 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recipient identification number cell (B2) to the new text value.
$ws.Range("B2").Value = "123456789"

# Move / record the active selection as B4 (was L10).
$ws.Range("B4").Select()
